$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row above the old row 6 (del/remove/destroy/remove row),
#    shifting rows 6-8 down to 7-9.
$ws.Rows("6").Insert()

# 2) Insert a new column before column G (the objectCache column), shifting it to H
#    and leaving a blank column F for the new "SF" backend.
$ws.Columns("F").Insert()

# 3) Header for the new "SF" column.
$ws.Range("F1").Value2 = "SF"
$ws.Range("F1").Style = $ws.Range("E1").Style

# 4) Fill column F (SF backend) with the same operations as column H (objectCache),
#    row by row, skipping the blank separator row (row 6).
$ws.Range("F2").Value2 = $ws.Range("H2").Value2
$ws.Range("F3").Value2 = $ws.Range("H3").Value2
$ws.Range("F4").Value2 = $ws.Range("H4").Value2
$ws.Range("F5").Value2 = $ws.Range("H5").Value2
$ws.Range("F7").Value2 = $ws.Range("H7").Value2
$ws.Range("F8").Value2 = $ws.Range("H8").Value2
$ws.Range("F9").Value2 = $ws.Range("H9").Value2

# 5) New "list_values" row between "list_keys" and "del", only present for the
#    objectCache column (H), plus two new trailing rows (relocate, about_me).
$ws.Range("H6").Value2 = "list_values"
$ws.Range("H10").Value2 = "relocate"
$ws.Range("H11").Value2 = "about_me"

# 6) Style the whole SF/objectCache detail range (italic, accent color) to match
#    the rest of the new column formatting, and extend the blank tail rows.
$italicRange = $ws.Range("F2:F9,H2:H17")
$italicRange.Font.Italic = $true
$italicRange.Font.ThemeColor = 10

# 7) Column widths: new column F, and column H keeps the old column G width.
$ws.Columns("F").ColumnWidth = 10.28515625
$ws.Columns("H").ColumnWidth = 11.85546875

# 8) Move the active selection the way the author left it.
$ws.Range("H11").Select()

Write-Host "done"
